$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the two price values in column D (D29, D30)
$ws.Range("D29").Value = 1001
$ws.Range("D30").Value = 1077
